$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must remain literal text (so numeric-/date-
# looking strings like "2025-01-01" or "7250585058" are not silently
# reinterpreted by Excel as a date serial / number). We briefly force a
# text number format so the auto-detection in the Value setter is bypassed,
# then restore the default "Normal" style so no stray style survives.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Remove the two extra student rows (rows 3 and 4), keeping only the header
# row and the single remaining student record (row 2).
$ws.Rows("3:4").Delete()

# Update the remaining student's record (row 2) with the new values.
Set-TextValue $ws.Range("B2") "2025-01-01"
$ws.Range("D2").Value = "testfather"
Set-TextValue $ws.Range("F2") "7250585058"
$ws.Range("H2").Value = "1,2"
Set-TextValue $ws.Range("K2") "150.00"
Set-TextValue $ws.Range("M2") "2025-03-05"
